$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the date/time values in column A (rows 2-25) forward by 77 days
# (44755 -> 44832), keeping the fractional time-of-day component intact.
for ($r = 2; $r -le 25; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 77
}

# Update the message-count values in column B for rows 10-13.
$ws.Cells.Item(10, 2).Value2 = 1
$ws.Cells.Item(11, 2).Value2 = 1
$ws.Cells.Item(12, 2).Value2 = 6
$ws.Cells.Item(13, 2).Value2 = 11
